# feat: add data validation to grade and student templates
#
# Mirrors the authored change to public/template/grade.xlsx:
#  - add a List data-validation ("1,2,S") to F2 (semesterSequence)
#  - correct the sample semesterSequence value in F2 from 2 to 1
#  - carry over the cosmetic re-save deltas that are reachable through the
#    Excel object model (zoom/selection, column widths, base font)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- primary feature: data validation on F2 (dropdown list 1 / 2 / S) -----
# Type:=xlValidateList(3), AlertStyle:=xlValidAlertStop(1), Operator left
# at its default (xlBetween) so no operator="" attribute is emitted for the
# list-type rule, Formula1 is the literal quoted in-cell list.
$ws.Range("F2").Validation.Add(3, 1, 1, '"1,2,S"')

# --- sample data fix that shipped alongside the validation ---------------
$ws.Range("F2").Value = 1

# --- cosmetic state that travelled with the same save ---------------------
# Window zoom 220% -> 175%, and the remembered selection moved to F5.
$excel.ActiveWindow.Zoom = 175
$ws.Range("F5").Select()

# Column widths nudged by the re-save (A: ~13.73 -> ~13.75, F: ~19.36 -> ~19.38 chars).
$ws.Columns("A").ColumnWidth = 12.83
$ws.Columns("F").ColumnWidth = 18.5

# Base font changed from "Aptos Narrow" to "Tahoma" for every styled cell
# already on the sheet (keeps the existing A1:F7 footprint untouched).
$cells = @("A1","B1","E1","F1","A2","B2","E2","F2","A3","B3","A4","B4","A5","B5","A6","B6","A7","B7")
foreach ($addr in $cells) {
  $ws.Range($addr).Font.Name = "Tahoma"
}
